$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.55469822883606
$ws.Range("B1").Value = 2.338839530944824
$ws.Range("C1").Value = 1.773891091346741
$ws.Range("D1").Value = 1.615631461143494
$ws.Range("E1").Value = 1.452238082885742
